$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.522.51"
$ws.Range("D3").Value = "3.667.29"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.01%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "624.21"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -7.31%  "
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "158.83"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  +0.03%  "
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.497"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("D13").Value = "4.286.87"
$ws.Range("E13").Value = "  -1.02%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "32.32"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "3.686.29"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "69.536.19"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("E18").Value = "  -0.11%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.89"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -2.77%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.35"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  +5.45%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "468.76"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("E22").Value = "  -0.68%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "79.64"
$cell.Style = $origStyle
$ws.Range("D24").Value = "3.814.27"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("E27").Value = "  -5.10%  "
$ws.Range("E28").Value = "  -6.05%  "
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("E30").Value = "  -4.45%  "
$ws.Range("E31").Value = "  +0.11%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.62"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("D35").Value = "3.672.74"
$ws.Range("E36").Value = "  -4.17%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.27"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  -3.34%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "178.47"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +2.62%  "
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -0.04%  "
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.22"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -1.77%  "
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.79"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -5.40%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0892"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  -2.46%  "
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.922"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "46.90"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "29.21"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +4.86%  "
$ws.Range("E47").Value = "  -2.45%  "
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.86"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  -4.99%  "
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.03"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  -5.37%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.259"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -3.26%  "
